$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (QH0N67zQ / Dep. Tachira vs Zamora) data has moved up into row 13
# with several odds values updated; the old row 14 is removed entirely.
$ws.Rows(14).Delete() | Out-Null

# Apply the updated odds/values for the affected rows (2,3,4,5,7,8,10)
# and the fully refreshed row 13 (id/time/teams + all odds columns).

# Row 2
$ws.Range("G2").Value = 2.1
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 2
$ws.Range("X2").Value = 9.5
$ws.Range("AD2").Value = 6
$ws.Range("AH2").Value = 9
$ws.Range("AP2").Value = 26
$ws.Range("AZ2").Value = 67

# Row 3
$ws.Range("G3").Value = 1.45
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 2.05
$ws.Range("L3").Value = 8
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 5
$ws.Range("Z3").Value = 9
$ws.Range("AB3").Value = 41
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 41
$ws.Range("AJ3").Value = 23
$ws.Range("AK3").Value = 101
$ws.Range("AL3").Value = 67
$ws.Range("AN3").Value = 3.2
$ws.Range("AO3").Value = 7.5
$ws.Range("AQ3").Value = 23
$ws.Range("AS3").Value = 251
$ws.Range("AU3").Value = 11
$ws.Range("AW3").Value = 8.5
$ws.Range("AY3").Value = 51
$ws.Range("AZ3").Value = 201
$ws.Range("BA3").Value = 251

# Row 4
$ws.Range("G4").Value = 3.8
$ws.Range("I4").Value = 1.95
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.63
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.8
$ws.Range("X4").Value = 19
$ws.Range("AC4").Value = 9
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 6.5
$ws.Range("AP4").Value = 34
$ws.Range("AQ4").Value = 81
$ws.Range("BA4").Value = 51

# Row 5
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75

# Row 7
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5

# Row 8
$ws.Range("J8").Value = 2.1

# Row 10
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 19
$ws.Range("Q10").Value = 1.53
$ws.Range("R10").Value = 2.4

# Row 13
$ws.Range("A13").Value = "QH0N67zQ"
$ws.Range("C13").Value = "20:00"
$ws.Range("E13").Value = "Dep. Tachira"
$ws.Range("F13").Value = "Zamora"
$ws.Range("G13").Value = 1.55
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 6.5
$ws.Range("J13").Value = 2.1
$ws.Range("K13").Value = 2.07
$ws.Range("L13").Value = 6.1
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 6.65
$ws.Range("O13").Value = 1.35
$ws.Range("P13").Value = 2.72
$ws.Range("Q13").Value = 2.02
$ws.Range("R13").Value = 1.62
$ws.Range("S13").Value = 1.4
$ws.Range("T13").Value = 2.5
$ws.Range("U13").Value = 2.02
$ws.Range("V13").Value = 1.62
$ws.Range("W13").Value = 5.4
$ws.Range("X13").Value = 6.4
$ws.Range("Y13").Value = 8.25
$ws.Range("Z13").Value = 11
$ws.Range("AA13").Value = 14
$ws.Range("AC13").Value = 7.9
$ws.Range("AD13").Value = 6.9
$ws.Range("AE13").Value = 18.5
$ws.Range("AF13").Value = 110
$ws.Range("AG13").Value = 201
$ws.Range("AH13").Value = 15
$ws.Range("AI13").Value = 40
$ws.Range("AJ13").Value = 20
$ws.Range("AK13").Value = 175
$ws.Range("AL13").Value = 80
$ws.Range("AM13").Value = 75
$ws.Range("AN13").Value = 3.25
$ws.Range("AO13").Value = 7.4
$ws.Range("AP13").Value = 18
$ws.Range("AQ13").Value = 24
$ws.Range("AR13").Value = 60
$ws.Range("AT13").Value = 2.47
$ws.Range("AU13").Value = 7.7
$ws.Range("AV13").Value = 75
$ws.Range("AW13").Value = 7.6
$ws.Range("AX13").Value = 37
$ws.Range("AY13").Value = 40
$ws.Range("AZ13").Value = 250
$ws.Range("BA13").Value = 250
$ws.Range("BB13").Value = 500
